$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.830119371414185
$ws.Range("B1").Value = 2.626540899276733
$ws.Range("C1").Value = 1.701378345489502
$ws.Range("D1").Value = 1.454130887985229
$ws.Range("E1").Value = 1.408189296722412
